$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.501.17'
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").Value = '2.372.92'
$ws.Range("E3").Value = '  -5.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '504.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.553'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.02%  '
$ws.Range("D9").Value = '2.395.14'
$ws.Range("E9").Value = '  -4.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0963'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -10.38%  '
$ws.Range("D14").Value = '2.800.61'
$ws.Range("E14").Value = '  -4.90%  '
$ws.Range("D15").Value = '56.381.06'
$ws.Range("E15").Value = '  -3.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.17%  '
$ws.Range("E17").Value = '  -3.21%  '
$ws.Range("D18").Value = '2.388.95'
$ws.Range("E18").Value = '  -4.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '312.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("D26").Value = '2.495.95'
$ws.Range("E26").Value = '  -4.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.375'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.03%  '
$ws.Range("E28").Value = '  -5.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("E31").Value = '  -3.56%  '
$ws.Range("D32").Value = '0.0₃0712'
$ws.Range("E32").Value = '  -5.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("E34").Value = '  -5.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.785'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '131.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '255.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.570'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0900'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0488'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0207'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.41%  '
